# Applies the "COMM_MTH_HASH_KEY" schema-design edit:
#   - Inserts a new row of F:I info (COMM_MTH_HASH_KEY / VARCHAR(40) / <blank> /
#     "Hash key value of commited method") right under row 20 on the
#     "Analysis DB" sheet, which pushes the whole existing F:I (Column
#     name/Data type/<blank>/Description) block that used to start at
#     row 23 down by one row (through row 47).
#   - Widens column F and nudges the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analysis DB")

# ---------------------------------------------------------------------------
# 1) Shift the existing F:I (Column name | Data type | <blank> | Description)
#    blocks down by one row: row 46 -> 47, 45 -> 46, ... , 23 -> 24.
#    Copy (not Cut) so we go bottom-up and never clobber a source row before
#    it has been read.
# ---------------------------------------------------------------------------
for ($r = 46; $r -ge 23; $r--) {
    $src = $ws.Range("F" + $r + ":I" + $r)
    $dst = $ws.Range("F" + ($r + 1) + ":I" + ($r + 1))
    $src.Copy($dst)
}

# The old F23:I23 content has now been duplicated onto F24:I24, so the
# vacated row must be cleared out (new row 23 has no F:I cells at all).
$ws.Range("F23:I23").ClearContents()
$ws.Range("F23:I23").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Populate the brand-new F21:I21 block (style matches the plain
#    bordered "s=1" cells used throughout this table).
# ---------------------------------------------------------------------------
$styleSrc = $ws.Range("B21")

$f21 = $ws.Range("F21")
$styleSrc.Copy($f21)
$f21.Value2 = "COMM_MTH_HASH_KEY"

$g21 = $ws.Range("G21")
$styleSrc.Copy($g21)
$g21.Value2 = "VARCHAR(40)"

$h21 = $ws.Range("H21")
$styleSrc.Copy($h21)
$h21.Value2 = ""

$i21 = $ws.Range("I21")
$styleSrc.Copy($i21)
$i21.Value2 = "Hash key value of commited method"

# ---------------------------------------------------------------------------
# 3) Row heights: the F:I block carries its row's auto-computed height
#    along with it, so re-stamp the heights for the rows that now hold
#    taller (wrapped) content and clear the ones that no longer need it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(23).RowHeight = 17
$ws.Rows.Item(24).RowHeight = 18
$ws.Rows.Item(29).RowHeight = 17
$ws.Rows.Item(30).RowHeight = 18
$ws.Rows.Item(31).RowHeight = 17
$ws.Rows.Item(32).RowHeight = 51
$ws.Rows.Item(38).RowHeight = 17
$ws.Rows.Item(39).RowHeight = 34
$ws.Rows.Item(41).RowHeight = 17
$ws.Rows.Item(42).RowHeight = 18
$ws.Rows.Item(43).RowHeight = 17
$ws.Rows.Item(44).RowHeight = 51
$ws.Rows.Item(46).RowHeight = 17
$ws.Rows.Item(47).RowHeight = 51

# ---------------------------------------------------------------------------
# 4) Column F widens to fit the longer "COMM_MTH_HASH_KEY" header.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 21.2857142857

# ---------------------------------------------------------------------------
# 5) Sheet view: scroll down a little and move the selection to the newly
#    added F21 cell.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("F21").Select()
